$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore cell C10 value to 1 (was 18)
$ws.Range("C10").Value = 1
